# Update "想去人数" (F column) counts on the two sheets that list every
# exhibition/event row: "展览" (Worksheet 1) and "全部类型" (Worksheet 4).
# Same underlying events are duplicated on both sheets (rows differ slightly
# after row 19 because sheet4 has one extra row inserted earlier).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# Sheet 1 (展览) updates: row -> new value
$ws1.Range("F2").Value  = 1330
$ws1.Range("F4").Value  = 14617
$ws1.Range("F5").Value  = 17736
$ws1.Range("F7").Value  = 72
$ws1.Range("F17").Value = 154
$ws1.Range("F19").Value = 1336
$ws1.Range("F24").Value = 7307
$ws1.Range("F28").Value = 1177
$ws1.Range("F30").Value = 5870
$ws1.Range("F32").Value = 48
$ws1.Range("F33").Value = 144
$ws1.Range("F35").Value = 227

# Sheet 4 (全部类型) updates: row -> new value
$ws4.Range("F2").Value  = 1330
$ws4.Range("F4").Value  = 14617
$ws4.Range("F5").Value  = 17736
$ws4.Range("F7").Value  = 72
$ws4.Range("F17").Value = 154
$ws4.Range("F19").Value = 1336
$ws4.Range("F25").Value = 7307
$ws4.Range("F29").Value = 1177
$ws4.Range("F32").Value = 5870
$ws4.Range("F34").Value = 48
$ws4.Range("F35").Value = 144
$ws4.Range("F37").Value = 227
